# Adding new progress as of date 04-Nov-2025:
#  - "PERIOD TO EXPIRE" (column H) drops by one day for every training row.
#  - "LAST UPDATE" (column I) moves from 03-Nov-2025 to 04-Nov-2025.
# Only the "Training Dashboard" sheet (rows 3-35) is affected.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Training Dashboard")

$lastRow = $ws.Cells.Item($ws.Rows.Count, 8).End(-4162).Row  # xlUp, column H

# --- Column H: decrement "PERIOD TO EXPIRE" by 1 for every data row --------
for ($r = 3; $r -le $lastRow; $r++) {
    $hCell = $ws.Cells.Item($r, 8)
    if ($hCell.Value2 -ne $null -and $hCell.Value2 -ne "") {
        $hCell.Value2 = $hCell.Value2 - 1
    }
}

# --- Column I: refresh "LAST UPDATE" to 04-Nov-2025 -------------------------
# Plain string assignment gets auto-parsed by Excel into a date serial
# (changing the cell's number format/style). To keep the cell as plain text
# with its original style, write the literal via a formula and then paste it
# back in as a static value (the same effect as Excel's Paste Values).
$dateRange = $ws.Range($ws.Cells.Item(3, 9), $ws.Cells.Item($lastRow, 9))
$dateRange.Formula = '="04-Nov-2025"'
$dateRange.Copy()
$dateRange.PasteSpecial(-4163)  # xlPasteValues
$excel.CutCopyMode = $false
